# SIQ workbook update
# Adds a new Software Interface Questionnaire entry (PO3_DGW_SIQ_14) asked
# and proposed by Amr on 22/2/2020, about alarm behaviour across modes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIQ")

# New question row (row 17, directly below the last existing entry, row 16)
$ws.Range("B17").Value = "Amr"
$ws.Range("C17").Value = "PO3_DGW_SIQ_14"
$ws.Range("D17").Value = "What to do when alarm is reached when mode is in another mode "
$ws.Range("F17").Value = "Amr"
$ws.Range("G17").Value = "22/2/2020"
$ws.Range("H17").Value = "Not answered"

# Reflect where the author's cursor ended up after entering the new row
$null = $ws.Range("D18").Select()

Write-Output "Added SIQ row 17 (PO3_DGW_SIQ_14)"
